$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Find the nearest  station to sell at that contains black market.", $true, $false, $false, $false, $false, $true, 1, $false, "Find all systems that are within X lightyears of System Y.", 2)
